$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row right before the current row 474, shifting all
# subsequent rows (old 474..507) down by one (to 475..508).
$ws.Rows.Item(474).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Cells.Item(474, 1).Value = 8
$ws.Cells.Item(474, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(474, 3).Value = "Coquimbo"
$ws.Cells.Item(474, 4).Value = 45021
$ws.Cells.Item(474, 5).Value = 4
$ws.Cells.Item(474, 6).Value = 100114013
$ws.Cells.Item(474, 7).Value = "Zanahoria"
$ws.Cells.Item(474, 8).Value = "Sin especificar"
$ws.Cells.Item(474, 9).Value = "Primera"
$ws.Cells.Item(474, 10).Value = 520
$ws.Cells.Item(474, 11).Value = 5500
$ws.Cells.Item(474, 12).Value = 6000
$ws.Cells.Item(474, 13).Value = 5750
$ws.Cells.Item(474, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(474, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(474, 16).Value = 288
$ws.Cells.Item(474, 17).Value = 20
$ws.Cells.Item(474, 18).Value = "Hortaliza"
